$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-by-row updates of Price (D) and Volume(1h) (E) columns,
# plus a B/C/D/E swap for rows 48-49 (FTXToken <-> FraxShare).
# D-column values that parse as plain decimals are prefixed with a
# leading apostrophe so Excel stores them as text (matching the
# original inline-string cell type) instead of auto-converting to a number.

# Row 2
$ws.Range("D2").Value = "37.428.83"
$ws.Range("E2").Value = "  -0.81%  "
# Row 3
$ws.Range("D3").Value = "2.067.47"
$ws.Range("E3").Value = "  -0.36%  "
# Row 4
$ws.Range("E4").Value = "  +0.09%  "
# Row 5
$ws.Range("D5").Value = "'232.03"
$ws.Range("E5").Value = "  -0.48%  "
# Row 6
$ws.Range("E6").Value = "  +0.90%  "
# Row 7
$ws.Range("E7").Value = "  +0.03%  "
# Row 8
$ws.Range("D8").Value = "'57.12"
$ws.Range("E8").Value = "  -2.48%  "
# Row 9
$ws.Range("D9").Value = "'0.388"
$ws.Range("E9").Value = "  -1.32%  "
# Row 10
$ws.Range("D10").Value = "'0.0776"
$ws.Range("E10").Value = "  -0.92%  "
# Row 11
$ws.Range("E11").Value = "  +1.44%  "
# Row 12
$ws.Range("D12").Value = "'14.79"
$ws.Range("E12").Value = "  +0.32%  "
# Row 13
$ws.Range("D13").Value = "2.374.60"
$ws.Range("E13").Value = "  -0.18%  "
# Row 14
$ws.Range("D14").Value = "'20.79"
$ws.Range("E14").Value = "  -0.82%  "
# Row 15
$ws.Range("D15").Value = "'0.762"
$ws.Range("E15").Value = "  -1.65%  "
# Row 16
$ws.Range("D16").Value = "'5.30"
$ws.Range("E16").Value = "  -1.32%  "
# Row 17
$ws.Range("D17").Value = "2.068.80"
$ws.Range("E17").Value = "  +0.28%  "
# Row 18
$ws.Range("D18").Value = "37.362.50"
$ws.Range("E18").Value = "  -0.76%  "
# Row 19
$ws.Range("D19").Value = "'70.33"
$ws.Range("E19").Value = "  -1.05%  "
# Row 20
$ws.Range("D20").Value = "'5.95"
$ws.Range("E20").Value = "  -2.97%  "
# Row 21
$ws.Range("D21").Value = "0.0₃0825"
$ws.Range("E21").Value = "  -1.14%  "
# Row 22
$ws.Range("D22").Value = "'227.77"
$ws.Range("E22").Value = "  -0.22%  "
# Row 23
$ws.Range("E23").Value = "  +0.01%  "
# Row 24
$ws.Range("E24").Value = "  -0.69%  "
# Row 25
$ws.Range("E25").Value = "  -1.39%  "
# Row 26
$ws.Range("D26").Value = "'9.58"
$ws.Range("E26").Value = "  +6.17%  "
# Row 27
$ws.Range("D27").Value = "'169.86"
$ws.Range("E27").Value = "  -0.60%  "
# Row 28
$ws.Range("D28").Value = "'0.132"
$ws.Range("E28").Value = "  -3.81%  "
# Row 29
$ws.Range("D29").Value = "'19.39"
$ws.Range("E29").Value = "  -0.54%  "
# Row 30
$ws.Range("E30").Value = "  -1.56%  "
# Row 31
$ws.Range("D31").Value = "'0.122"
$ws.Range("E31").Value = "  +0.80%  "
# Row 32
$ws.Range("D32").Value = "'4.59"
$ws.Range("E32").Value = "  -2.04%  "
# Row 33
$ws.Range("D33").Value = "'0.0632"
$ws.Range("E33").Value = "  -0.08%  "
# Row 34
$ws.Range("E34").Value = "  -0.74%  "
# Row 35
$ws.Range("D35").Value = "'2.46"
$ws.Range("E35").Value = "  -1.43%  "
# Row 36
$ws.Range("D36").Value = "'1.82"
$ws.Range("E36").Value = "  -0.31%  "
# Row 37
$ws.Range("D37").Value = "'3.30"
$ws.Range("E37").Value = "  -2.55%  "
# Row 38
$ws.Range("E38").Value = "  -0.05%  "
# Row 39
$ws.Range("D39").Value = "'5.25"
$ws.Range("E39").Value = "  -1.24%  "
# Row 40
$ws.Range("D40").Value = "'0.0229"
$ws.Range("E40").Value = "  +6.46%  "
# Row 41
$ws.Range("D41").Value = "'99.51"
$ws.Range("E41").Value = "  -0.81%  "
# Row 42
$ws.Range("D42").Value = "'2.90"
$ws.Range("E42").Value = "  +0.74%  "
# Row 43
$ws.Range("D43").Value = "'1.19"
$ws.Range("E43").Value = "  +3.57%  "
# Row 44
$ws.Range("D44").Value = "'0.0947"
$ws.Range("E44").Value = "  -2.95%  "
# Row 45
$ws.Range("D45").Value = "1.475.74"
$ws.Range("E45").Value = "  +2.27%  "
# Row 46
$ws.Range("D46").Value = "'16.67"
$ws.Range("E46").Value = "  -0.37%  "
# Row 47
$ws.Range("E47").Value = "  -1.94%  "
# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.26"
$ws.Range("E48").Value = "  -2.26%  "
# Row 49
$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").Value = "'3.94"
$ws.Range("E49").Value = "  -6.51%  "
# Row 50
$ws.Range("D50").Value = "'2.93"
$ws.Range("E50").Value = "  -1.85%  "
# Row 51
$ws.Range("D51").Value = "2.259.07"
$ws.Range("E51").Value = "  -0.24%  "
